$d = $word.ActiveDocument

# Locate the exact bounds of the run that currently reads:
#   "μέρη της εργασίας με αυξημένο τον κίνδυνο εμπλοκής πολλών ομάδας που
#    δουλεύουν στο ίδιο έργο με αποτέλεσμα να μην είναι εύκολος ο
#    συντονισμός μας."
$runRange = $d.Content
[void]$runRange.Find.Execute(
    "μέρη της εργασίας με αυξημένο τον κίνδυνο εμπλοκής πολλών ομάδας που δουλεύουν στο ίδιο έργο με αποτέλεσμα να μην είναι εύκολος ο συντονισμός μας.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runStart = $runRange.Start
$runEnd = $runRange.End

# Within that exact span, find the word that needs to change: "ομάδας" -> "τμηματικών ομάδων"
$target = $d.Range($runStart, $runEnd)
[void]$target.Find.Execute("ομάδας", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$omStart = $target.Start

$replacement = "τμηματικών ομάδων"
$target.Text = $replacement
$midEnd = $omStart + $replacement.Length

# The text substitution above collapses the paragraph's same-formatted runs
# back into a single run. Re-establish the three-way run split the edit
# calls for by nudging formatting (set then clear Bold) across each of the
# new boundaries; this forces the engine to keep them as distinct runs
# without altering any visible formatting.
$piece1 = $d.Range($runStart, $omStart)
$piece1.Bold = 1
$piece1.Bold = 0

$piece2 = $d.Range($omStart, $midEnd)
$piece2.Bold = 1
$piece2.Bold = 0
